$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9142543333333334
$ws.Range("H2").Value = 2.742763
$ws.Range("I2").Value = 0.1175834869881751
$ws.Range("J2").Value = 0.1175834869881751
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.604331
$ws.Range("N2").Value = 40.812993
$ws.Range("O2").Value = 0.8107276168878804
$ws.Range("P2").Value = 0.8107276168878805
$ws.Range("Q2").Value = 12.437818568851
$ws.Range("R2").Value = 111.940367119659
$ws.Range("S2").Value = 0.09532818019129027
$ws.Range("T2").Value = 0.09532818019129027

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9142543333333334
$ws.Range("H3").Value = 2.742763
$ws.Range("I3").Value = 0.1175834869881751
$ws.Range("J3").Value = 0.1175834869881751
$ws.Range("O3").Value = 0.06327311690486458
$ws.Range("P3").Value = 0.06327311690486459
$ws.Range("Q3").Value = 0.970707709908
$ws.Range("R3").Value = 8.736369389172
$ws.Range("S3").Value = 0.007439873718284424
$ws.Range("T3").Value = 0.007439873718284425

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9142543333333334
$ws.Range("H4").Value = 2.742763
$ws.Range("I4").Value = 0.1175834869881751
$ws.Range("J4").Value = 0.1175834869881751
$ws.Range("M4").Value = 1.995771333333333
$ws.Range("N4").Value = 5.987314
$ws.Range("O4").Value = 0.1189346934389115
$ws.Range("P4").Value = 0.1189346934389116
$ws.Range("Q4").Value = 1.824642589842444
$ws.Range("R4").Value = 16.421783308582
$ws.Range("S4").Value = 0.01398475597841684
$ws.Range("T4").Value = 0.01398475597841685

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9142543333333334
$ws.Range("H5").Value = 2.742763
$ws.Range("I5").Value = 0.1175834869881751
$ws.Range("J5").Value = 0.1175834869881751
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1185463333333334
$ws.Range("N5").Value = 0.355639
$ws.Range("O5").Value = 0.007064572768343379
$ws.Range("P5").Value = 0.007064572768343379
$ws.Range("Q5").Value = 0.1083814989507778
$ws.Range("R5").Value = 0.9754334905570001
$ws.Range("S5").Value = 0.0008306771001835196
$ws.Range("T5").Value = 0.0008306771001835196

# Row 6
$ws.Range("I6").Value = 0.6206849497708361
$ws.Range("J6").Value = 0.620684949770836
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 13.604331
$ws.Range("N6").Value = 40.812993
$ws.Range("O6").Value = 0.8107276168878804
$ws.Range("P6").Value = 0.8107276168878805
$ws.Range("Q6").Value = 65.65519522688099
$ws.Range("R6").Value = 590.8967570419289
$ws.Range("S6").Value = 0.5032064301658837
$ws.Range("T6").Value = 0.5032064301658837

# Row 7
$ws.Range("I7").Value = 0.6206849497708361
$ws.Range("J7").Value = 0.620684949770836
$ws.Range("O7").Value = 0.06327311690486458
$ws.Range("P7").Value = 0.06327311690486459
$ws.Range("S7").Value = 0.03927267138794011
$ws.Range("T7").Value = 0.03927267138794011

# Row 8
$ws.Range("I8").Value = 0.6206849497708361
$ws.Range("J8").Value = 0.620684949770836
$ws.Range("M8").Value = 1.995771333333333
$ws.Range("N8").Value = 5.987314
$ws.Range("O8").Value = 0.1189346934389115
$ws.Range("P8").Value = 0.1189346934389116
$ws.Range("Q8").Value = 9.631694239004664
$ws.Range("R8").Value = 86.68524815104199
$ws.Range("S8").Value = 0.07382097422314059
$ws.Range("T8").Value = 0.0738209742231406

# Row 9
$ws.Range("I9").Value = 0.6206849497708361
$ws.Range("J9").Value = 0.620684949770836
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1185463333333334
$ws.Range("N9").Value = 0.355639
$ws.Range("O9").Value = 0.007064572768343379
$ws.Range("P9").Value = 0.007064572768343379
$ws.Range("Q9").Value = 0.5721106505296667
$ws.Range("R9").Value = 5.148995854767001
$ws.Range("S9").Value = 0.004384873993871627
$ws.Range("T9").Value = 0.004384873993871626

# Row 10
$ws.Range("G10").Value = 1.986145
$ws.Range("H10").Value = 5.958435
$ws.Range("I10").Value = 0.2554407961214246
$ws.Range("J10").Value = 0.2554407961214246
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.604331
$ws.Range("N10").Value = 40.812993
$ws.Range("O10").Value = 0.8107276168878804
$ws.Range("P10").Value = 0.8107276168878805
$ws.Range("Q10").Value = 27.020173993995
$ws.Range("R10").Value = 243.181565945955
$ws.Range("S10").Value = 0.2070929078954655
$ws.Range("T10").Value = 0.2070929078954655

# Row 11
$ws.Range("G11").Value = 1.986145
$ws.Range("H11").Value = 5.958435
$ws.Range("I11").Value = 0.2554407961214246
$ws.Range("J11").Value = 0.2554407961214246
$ws.Range("O11").Value = 0.06327311690486458
$ws.Range("P11").Value = 0.06327311690486459
$ws.Range("Q11").Value = 2.10878548146
$ws.Range("R11").Value = 18.97906933314
$ws.Range("S11").Value = 0.01616253535526258
$ws.Range("T11").Value = 0.01616253535526258

# Row 12
$ws.Range("G12").Value = 1.986145
$ws.Range("H12").Value = 5.958435
$ws.Range("I12").Value = 0.2554407961214246
$ws.Range("J12").Value = 0.2554407961214246
$ws.Range("M12").Value = 1.995771333333333
$ws.Range("N12").Value = 5.987314
$ws.Range("O12").Value = 0.1189346934389115
$ws.Range("P12").Value = 0.1189346934389116
$ws.Range("Q12").Value = 3.963891254843332
$ws.Range("R12").Value = 35.67502129358999
$ws.Range("S12").Value = 0.03038077277849313
$ws.Range("T12").Value = 0.03038077277849314

# Row 13
$ws.Range("G13").Value = 1.986145
$ws.Range("H13").Value = 5.958435
$ws.Range("I13").Value = 0.2554407961214246
$ws.Range("J13").Value = 0.2554407961214246
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1185463333333334
$ws.Range("N13").Value = 0.355639
$ws.Range("O13").Value = 0.007064572768343379
$ws.Range("P13").Value = 0.007064572768343379
$ws.Range("Q13").Value = 0.2354502072183333
$ws.Range("R13").Value = 2.119051864965
$ws.Range("S13").Value = 0.001804580092203369
$ws.Range("T13").Value = 0.001804580092203369

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.04891300000000001
$ws.Range("H14").Value = 0.146739
$ws.Range("I14").Value = 0.006290767119564404
$ws.Range("J14").Value = 0.006290767119564403
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 13.604331
$ws.Range("N14").Value = 40.812993
$ws.Range("O14").Value = 0.8107276168878804
$ws.Range("P14").Value = 0.8107276168878805
$ws.Range("Q14").Value = 0.6654286422030001
$ws.Range("R14").Value = 5.988857779827
$ws.Range("S14").Value = 0.005100098635241085
$ws.Range("T14").Value = 0.005100098635241085

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.04891300000000001
$ws.Range("H15").Value = 0.146739
$ws.Range("I15").Value = 0.006290767119564404
$ws.Range("J15").Value = 0.006290767119564403
$ws.Range("O15").Value = 0.06327311690486458
$ws.Range("P15").Value = 0.06327311690486459
$ws.Range("Q15").Value = 0.051933279924
$ws.Range("R15").Value = 0.467399519316
$ws.Range("S15").Value = 0.0003980364433774767
$ws.Range("T15").Value = 0.0003980364433774767

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.04891300000000001
$ws.Range("H16").Value = 0.146739
$ws.Range("I16").Value = 0.006290767119564404
$ws.Range("J16").Value = 0.006290767119564403
$ws.Range("M16").Value = 1.995771333333333
$ws.Range("N16").Value = 5.987314
$ws.Range("O16").Value = 0.1189346934389115
$ws.Range("P16").Value = 0.1189346934389116
$ws.Range("Q16").Value = 0.09761916322733333
$ws.Range("R16").Value = 0.878572469046
$ws.Range("S16").Value = 0.0007481904588609768
$ws.Range("T16").Value = 0.0007481904588609769

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.04891300000000001
$ws.Range("H17").Value = 0.146739
$ws.Range("I17").Value = 0.006290767119564404
$ws.Range("J17").Value = 0.006290767119564403
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.1185463333333334
$ws.Range("N17").Value = 0.355639
$ws.Range("O17").Value = 0.007064572768343379
$ws.Range("P17").Value = 0.007064572768343379
$ws.Range("Q17").Value = 0.005798456802333335
$ws.Range("R17").Value = 0.05218611122100001
$ws.Range("S17").Value = 0.00004444158208486461
$ws.Range("T17").Value = 0.0000444415820848646

